$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 47,20

$data[0,0] = 5
$data[0,1] = 'Macroferia Regional de Talca'
$data[0,2] = 'Maule'
$data[0,3] = 44292
$data[0,4] = 7
$data[0,5] = 'Fruta'
$data[0,6] = 100107
$data[0,7] = 'Otros'
$data[0,8] = 100107011
$data[0,9] = 'Tuna'
$data[0,10] = 'Sin especificar'
$data[0,11] = 'Especial'
$data[0,12] = 150
$data[0,13] = 16000
$data[0,14] = 16000
$data[0,15] = 16000
$data[0,16] = '$/caja 18 kilos'
$data[0,17] = 'Provincia de Melipilla'
$data[0,18] = 889
$data[0,19] = 18

$data[1,0] = 5
$data[1,1] = 'Macroferia Regional de Talca'
$data[1,2] = 'Maule'
$data[1,3] = 44292
$data[1,4] = 7
$data[1,5] = 'Fruta'
$data[1,6] = 100107
$data[1,7] = 'Otros'
$data[1,8] = 100107011
$data[1,9] = 'Tuna'
$data[1,10] = 'Sin especificar'
$data[1,11] = 'Primera'
$data[1,12] = 80
$data[1,13] = 14000
$data[1,14] = 14000
$data[1,15] = 14000
$data[1,16] = '$/caja 18 kilos'
$data[1,17] = 'Provincia de Melipilla'
$data[1,18] = 778
$data[1,19] = 18

$data[2,0] = 5
$data[2,1] = 'Macroferia Regional de Talca'
$data[2,2] = 'Maule'
$data[2,3] = 44252
$data[2,4] = 7
$data[2,5] = 'Fruta'
$data[2,6] = 100107
$data[2,7] = 'Otros'
$data[2,8] = 100107011
$data[2,9] = 'Tuna'
$data[2,10] = 'Sin especificar'
$data[2,11] = 'Primera'
$data[2,12] = 140
$data[2,13] = 13000
$data[2,14] = 13000
$data[2,15] = 13000
$data[2,16] = '$/caja 18 kilos'
$data[2,17] = 'Provincia de Melipilla'
$data[2,18] = 722
$data[2,19] = 18

$data[3,0] = 5
$data[3,1] = 'Macroferia Regional de Talca'
$data[3,2] = 'Maule'
$data[3,3] = 44277
$data[3,4] = 7
$data[3,5] = 'Fruta'
$data[3,6] = 100107
$data[3,7] = 'Otros'
$data[3,8] = 100107011
$data[3,9] = 'Tuna'
$data[3,10] = 'Sin especificar'
$data[3,11] = 'Especial'
$data[3,12] = 200
$data[3,13] = 15000
$data[3,14] = 15000
$data[3,15] = 15000
$data[3,16] = '$/caja 18 kilos'
$data[3,17] = 'Provincia de Limarí'
$data[3,18] = 833
$data[3,19] = 18

$data[4,0] = 5
$data[4,1] = 'Macroferia Regional de Talca'
$data[4,2] = 'Maule'
$data[4,3] = 44299
$data[4,4] = 7
$data[4,5] = 'Fruta'
$data[4,6] = 100107
$data[4,7] = 'Otros'
$data[4,8] = 100107011
$data[4,9] = 'Tuna'
$data[4,10] = 'Sin especificar'
$data[4,11] = 'Especial'
$data[4,12] = 170
$data[4,13] = 18000
$data[4,14] = 18000
$data[4,15] = 18000
$data[4,16] = '$/caja 18 kilos'
$data[4,17] = 'Provincia de Melipilla'
$data[4,18] = 1000
$data[4,19] = 18

$data[5,0] = 5
$data[5,1] = 'Macroferia Regional de Talca'
$data[5,2] = 'Maule'
$data[5,3] = 44299
$data[5,4] = 7
$data[5,5] = 'Fruta'
$data[5,6] = 100107
$data[5,7] = 'Otros'
$data[5,8] = 100107011
$data[5,9] = 'Tuna'
$data[5,10] = 'Sin especificar'
$data[5,11] = 'Primera'
$data[5,12] = 100
$data[5,13] = 16000
$data[5,14] = 16000
$data[5,15] = 16000
$data[5,16] = '$/caja 18 kilos'
$data[5,17] = 'Provincia de Melipilla'
$data[5,18] = 889
$data[5,19] = 18

$data[6,0] = 5
$data[6,1] = 'Macroferia Regional de Talca'
$data[6,2] = 'Maule'
$data[6,3] = 44279
$data[6,4] = 7
$data[6,5] = 'Fruta'
$data[6,6] = 100107
$data[6,7] = 'Otros'
$data[6,8] = 100107011
$data[6,9] = 'Tuna'
$data[6,10] = 'Sin especificar'
$data[6,11] = 'Especial'
$data[6,12] = 50
$data[6,13] = 14000
$data[6,14] = 14000
$data[6,15] = 14000
$data[6,16] = '$/caja 18 kilos'
$data[6,17] = 'Provincia de Melipilla'
$data[6,18] = 778
$data[6,19] = 18

$data[7,0] = 5
$data[7,1] = 'Macroferia Regional de Talca'
$data[7,2] = 'Maule'
$data[7,3] = 44279
$data[7,4] = 7
$data[7,5] = 'Fruta'
$data[7,6] = 100107
$data[7,7] = 'Otros'
$data[7,8] = 100107011
$data[7,9] = 'Tuna'
$data[7,10] = 'Sin especificar'
$data[7,11] = 'Primera'
$data[7,12] = 100
$data[7,13] = 12000
$data[7,14] = 12000
$data[7,15] = 12000
$data[7,16] = '$/caja 18 kilos'
$data[7,17] = 'Provincia de Melipilla'
$data[7,18] = 667
$data[7,19] = 18

$data[8,0] = 5
$data[8,1] = 'Macroferia Regional de Talca'
$data[8,2] = 'Maule'
$data[8,3] = 44222
$data[8,4] = 7
$data[8,5] = 'Fruta'
$data[8,6] = 100107
$data[8,7] = 'Otros'
$data[8,8] = 100107011
$data[8,9] = 'Tuna'
$data[8,10] = 'Sin especificar'
$data[8,11] = 'Primera'
$data[8,12] = 100
$data[8,13] = 18000
$data[8,14] = 18000
$data[8,15] = 18000
$data[8,16] = '$/caja 16 kilos'
$data[8,17] = 'Provincia de Limarí'
$data[8,18] = 1125
$data[8,19] = 16

$data[9,0] = 5
$data[9,1] = 'Macroferia Regional de Talca'
$data[9,2] = 'Maule'
$data[9,3] = 44309
$data[9,4] = 7
$data[9,5] = 'Fruta'
$data[9,6] = 100107
$data[9,7] = 'Otros'
$data[9,8] = 100107011
$data[9,9] = 'Tuna'
$data[9,10] = 'Sin especificar'
$data[9,11] = 'Especial'
$data[9,12] = 100
$data[9,13] = 20000
$data[9,14] = 20000
$data[9,15] = 20000
$data[9,16] = '$/caja 18 kilos'
$data[9,17] = 'Provincia de Melipilla'
$data[9,18] = 1111
$data[9,19] = 18

$data[10,0] = 5
$data[10,1] = 'Macroferia Regional de Talca'
$data[10,2] = 'Maule'
$data[10,3] = 44309
$data[10,4] = 7
$data[10,5] = 'Fruta'
$data[10,6] = 100107
$data[10,7] = 'Otros'
$data[10,8] = 100107011
$data[10,9] = 'Tuna'
$data[10,10] = 'Sin especificar'
$data[10,11] = 'Primera'
$data[10,12] = 60
$data[10,13] = 18000
$data[10,14] = 18000
$data[10,15] = 18000
$data[10,16] = '$/caja 18 kilos'
$data[10,17] = 'Provincia de Melipilla'
$data[10,18] = 1000
$data[10,19] = 18

$data[11,0] = 5
$data[11,1] = 'Macroferia Regional de Talca'
$data[11,2] = 'Maule'
$data[11,3] = 44300
$data[11,4] = 7
$data[11,5] = 'Fruta'
$data[11,6] = 100107
$data[11,7] = 'Otros'
$data[11,8] = 100107011
$data[11,9] = 'Tuna'
$data[11,10] = 'Sin especificar'
$data[11,11] = 'Especial'
$data[11,12] = 120
$data[11,13] = 18000
$data[11,14] = 18000
$data[11,15] = 18000
$data[11,16] = '$/caja 18 kilos'
$data[11,17] = 'Provincia de Melipilla'
$data[11,18] = 1000
$data[11,19] = 18

$data[12,0] = 5
$data[12,1] = 'Macroferia Regional de Talca'
$data[12,2] = 'Maule'
$data[12,3] = 44300
$data[12,4] = 7
$data[12,5] = 'Fruta'
$data[12,6] = 100107
$data[12,7] = 'Otros'
$data[12,8] = 100107011
$data[12,9] = 'Tuna'
$data[12,10] = 'Sin especificar'
$data[12,11] = 'Primera'
$data[12,12] = 100
$data[12,13] = 16000
$data[12,14] = 16000
$data[12,15] = 16000
$data[12,16] = '$/caja 18 kilos'
$data[12,17] = 'Provincia de Melipilla'
$data[12,18] = 889
$data[12,19] = 18

$data[13,0] = 5
$data[13,1] = 'Macroferia Regional de Talca'
$data[13,2] = 'Maule'
$data[13,3] = 44273
$data[13,4] = 7
$data[13,5] = 'Fruta'
$data[13,6] = 100107
$data[13,7] = 'Otros'
$data[13,8] = 100107011
$data[13,9] = 'Tuna'
$data[13,10] = 'Sin especificar'
$data[13,11] = 'Especial'
$data[13,12] = 40
$data[13,13] = 15000
$data[13,14] = 15000
$data[13,15] = 15000
$data[13,16] = '$/caja 16 kilos'
$data[13,17] = 'Provincia de Melipilla'
$data[13,18] = 938
$data[13,19] = 16

$data[14,0] = 5
$data[14,1] = 'Macroferia Regional de Talca'
$data[14,2] = 'Maule'
$data[14,3] = 44273
$data[14,4] = 7
$data[14,5] = 'Fruta'
$data[14,6] = 100107
$data[14,7] = 'Otros'
$data[14,8] = 100107011
$data[14,9] = 'Tuna'
$data[14,10] = 'Sin especificar'
$data[14,11] = 'Primera'
$data[14,12] = 50
$data[14,13] = 13000
$data[14,14] = 13000
$data[14,15] = 13000
$data[14,16] = '$/caja 16 kilos'
$data[14,17] = 'Provincia de Melipilla'
$data[14,18] = 812
$data[14,19] = 16

$data[15,0] = 5
$data[15,1] = 'Macroferia Regional de Talca'
$data[15,2] = 'Maule'
$data[15,3] = 44273
$data[15,4] = 7
$data[15,5] = 'Fruta'
$data[15,6] = 100107
$data[15,7] = 'Otros'
$data[15,8] = 100107011
$data[15,9] = 'Tuna'
$data[15,10] = 'Sin especificar'
$data[15,11] = 'Segunda'
$data[15,12] = 60
$data[15,13] = 10000
$data[15,14] = 10000
$data[15,15] = 10000
$data[15,16] = '$/caja 16 kilos'
$data[15,17] = 'Provincia de Melipilla'
$data[15,18] = 625
$data[15,19] = 16

$data[16,0] = 5
$data[16,1] = 'Macroferia Regional de Talca'
$data[16,2] = 'Maule'
$data[16,3] = 45001
$data[16,4] = 7
$data[16,5] = 'Fruta'
$data[16,6] = 100107
$data[16,7] = 'Otros'
$data[16,8] = 100107011
$data[16,9] = 'Tuna'
$data[16,10] = 'Sin especificar'
$data[16,11] = 'Especial'
$data[16,12] = 150
$data[16,13] = 13000
$data[16,14] = 13000
$data[16,15] = 13000
$data[16,16] = '$/caja 18 kilos'
$data[16,17] = 'Provincia de Melipilla'
$data[16,18] = 722
$data[16,19] = 18

$data[17,0] = 5
$data[17,1] = 'Macroferia Regional de Talca'
$data[17,2] = 'Maule'
$data[17,3] = 45001
$data[17,4] = 7
$data[17,5] = 'Fruta'
$data[17,6] = 100107
$data[17,7] = 'Otros'
$data[17,8] = 100107011
$data[17,9] = 'Tuna'
$data[17,10] = 'Sin especificar'
$data[17,11] = 'Primera'
$data[17,12] = 100
$data[17,13] = 11000
$data[17,14] = 11000
$data[17,15] = 11000
$data[17,16] = '$/caja 18 kilos'
$data[17,17] = 'Provincia de Melipilla'
$data[17,18] = 611
$data[17,19] = 18

$data[18,0] = 5
$data[18,1] = 'Macroferia Regional de Talca'
$data[18,2] = 'Maule'
$data[18,3] = 44630
$data[18,4] = 7
$data[18,5] = 'Fruta'
$data[18,6] = 100107
$data[18,7] = 'Otros'
$data[18,8] = 100107011
$data[18,9] = 'Tuna'
$data[18,10] = 'Sin especificar'
$data[18,11] = 'Especial'
$data[18,12] = 150
$data[18,13] = 20000
$data[18,14] = 20000
$data[18,15] = 20000
$data[18,16] = '$/caja 20 kilos'
$data[18,17] = 'Provincia de Limarí'
$data[18,18] = 1000
$data[18,19] = 20

$data[19,0] = 5
$data[19,1] = 'Macroferia Regional de Talca'
$data[19,2] = 'Maule'
$data[19,3] = 45089
$data[19,4] = 7
$data[19,5] = 'Fruta'
$data[19,6] = 100107
$data[19,7] = 'Otros'
$data[19,8] = 100107011
$data[19,9] = 'Tuna'
$data[19,10] = 'Sin especificar'
$data[19,11] = 'Primera'
$data[19,12] = 30
$data[19,13] = 22000
$data[19,14] = 22000
$data[19,15] = 22000
$data[19,16] = '$/caja 18 kilos'
$data[19,17] = 'Provincia de Melipilla'
$data[19,18] = 1222
$data[19,19] = 18

$data[20,0] = 5
$data[20,1] = 'Macroferia Regional de Talca'
$data[20,2] = 'Maule'
$data[20,3] = 44330
$data[20,4] = 7
$data[20,5] = 'Fruta'
$data[20,6] = 100107
$data[20,7] = 'Otros'
$data[20,8] = 100107011
$data[20,9] = 'Tuna'
$data[20,10] = 'Sin especificar'
$data[20,11] = 'Primera'
$data[20,12] = 50
$data[20,13] = 23000
$data[20,14] = 23000
$data[20,15] = 23000
$data[20,16] = '$/caja 18 kilos'
$data[20,17] = 'Provincia de Melipilla'
$data[20,18] = 1278
$data[20,19] = 18

$data[21,0] = 5
$data[21,1] = 'Macroferia Regional de Talca'
$data[21,2] = 'Maule'
$data[21,3] = 44291
$data[21,4] = 7
$data[21,5] = 'Fruta'
$data[21,6] = 100107
$data[21,7] = 'Otros'
$data[21,8] = 100107011
$data[21,9] = 'Tuna'
$data[21,10] = 'Sin especificar'
$data[21,11] = 'Extra (doble especial)'
$data[21,12] = 250
$data[21,13] = 18000
$data[21,14] = 18000
$data[21,15] = 18000
$data[21,16] = '$/caja 18 kilos'
$data[21,17] = 'Provincia de Melipilla'
$data[21,18] = 1000
$data[21,19] = 18

$data[22,0] = 5
$data[22,1] = 'Macroferia Regional de Talca'
$data[22,2] = 'Maule'
$data[22,3] = 45093
$data[22,4] = 7
$data[22,5] = 'Fruta'
$data[22,6] = 100107
$data[22,7] = 'Otros'
$data[22,8] = 100107011
$data[22,9] = 'Tuna'
$data[22,10] = 'Sin especificar'
$data[22,11] = 'Primera'
$data[22,12] = 40
$data[22,13] = 22000
$data[22,14] = 22000
$data[22,15] = 22000
$data[22,16] = '$/caja 18 kilos'
$data[22,17] = 'Provincia de Melipilla'
$data[22,18] = 1222
$data[22,19] = 18

$data[23,0] = 5
$data[23,1] = 'Macroferia Regional de Talca'
$data[23,2] = 'Maule'
$data[23,3] = 44258
$data[23,4] = 7
$data[23,5] = 'Fruta'
$data[23,6] = 100107
$data[23,7] = 'Otros'
$data[23,8] = 100107011
$data[23,9] = 'Tuna'
$data[23,10] = 'Sin especificar'
$data[23,11] = 'Primera'
$data[23,12] = 100
$data[23,13] = 14000
$data[23,14] = 14000
$data[23,15] = 14000
$data[23,16] = '$/caja 18 kilos'
$data[23,17] = 'Provincia de Limarí'
$data[23,18] = 778
$data[23,19] = 18

$data[24,0] = 5
$data[24,1] = 'Macroferia Regional de Talca'
$data[24,2] = 'Maule'
$data[24,3] = 44315
$data[24,4] = 7
$data[24,5] = 'Fruta'
$data[24,6] = 100107
$data[24,7] = 'Otros'
$data[24,8] = 100107011
$data[24,9] = 'Tuna'
$data[24,10] = 'Sin especificar'
$data[24,11] = 'Especial'
$data[24,12] = 50
$data[24,13] = 24000
$data[24,14] = 24000
$data[24,15] = 24000
$data[24,16] = '$/caja 18 kilos'
$data[24,17] = 'Provincia de Melipilla'
$data[24,18] = 1333
$data[24,19] = 18

$data[25,0] = 5
$data[25,1] = 'Macroferia Regional de Talca'
$data[25,2] = 'Maule'
$data[25,3] = 44315
$data[25,4] = 7
$data[25,5] = 'Fruta'
$data[25,6] = 100107
$data[25,7] = 'Otros'
$data[25,8] = 100107011
$data[25,9] = 'Tuna'
$data[25,10] = 'Sin especificar'
$data[25,11] = 'Primera'
$data[25,12] = 50
$data[25,13] = 20000
$data[25,14] = 20000
$data[25,15] = 20000
$data[25,16] = '$/caja 18 kilos'
$data[25,17] = 'Provincia de Melipilla'
$data[25,18] = 1111
$data[25,19] = 18

$data[26,0] = 5
$data[26,1] = 'Macroferia Regional de Talca'
$data[26,2] = 'Maule'
$data[26,3] = 44274
$data[26,4] = 7
$data[26,5] = 'Fruta'
$data[26,6] = 100107
$data[26,7] = 'Otros'
$data[26,8] = 100107011
$data[26,9] = 'Tuna'
$data[26,10] = 'Sin especificar'
$data[26,11] = 'Especial'
$data[26,12] = 200
$data[26,13] = 14000
$data[26,14] = 14000
$data[26,15] = 14000
$data[26,16] = '$/caja 16 kilos'
$data[26,17] = 'Provincia de Melipilla'
$data[26,18] = 875
$data[26,19] = 16

$data[27,0] = 5
$data[27,1] = 'Macroferia Regional de Talca'
$data[27,2] = 'Maule'
$data[27,3] = 44274
$data[27,4] = 7
$data[27,5] = 'Fruta'
$data[27,6] = 100107
$data[27,7] = 'Otros'
$data[27,8] = 100107011
$data[27,9] = 'Tuna'
$data[27,10] = 'Sin especificar'
$data[27,11] = 'Primera'
$data[27,12] = 130
$data[27,13] = 12000
$data[27,14] = 12000
$data[27,15] = 12000
$data[27,16] = '$/caja 16 kilos'
$data[27,17] = 'Provincia de Melipilla'
$data[27,18] = 750
$data[27,19] = 16

$data[28,0] = 5
$data[28,1] = 'Macroferia Regional de Talca'
$data[28,2] = 'Maule'
$data[28,3] = 44699
$data[28,4] = 7
$data[28,5] = 'Fruta'
$data[28,6] = 100107
$data[28,7] = 'Otros'
$data[28,8] = 100107011
$data[28,9] = 'Tuna'
$data[28,10] = 'Sin especificar'
$data[28,11] = 'Especial'
$data[28,12] = 150
$data[28,13] = 22000
$data[28,14] = 22000
$data[28,15] = 22000
$data[28,16] = '$/caja 18 kilos'
$data[28,17] = 'Provincia de Limarí'
$data[28,18] = 1222
$data[28,19] = 18

$data[29,0] = 5
$data[29,1] = 'Macroferia Regional de Talca'
$data[29,2] = 'Maule'
$data[29,3] = 44645
$data[29,4] = 7
$data[29,5] = 'Fruta'
$data[29,6] = 100107
$data[29,7] = 'Otros'
$data[29,8] = 100107011
$data[29,9] = 'Tuna'
$data[29,10] = 'Sin especificar'
$data[29,11] = 'Primera'
$data[29,12] = 200
$data[29,13] = 16000
$data[29,14] = 16000
$data[29,15] = 16000
$data[29,16] = '$/caja 18 kilos'
$data[29,17] = 'Provincia de Limarí'
$data[29,18] = 889
$data[29,19] = 18

$data[30,0] = 5
$data[30,1] = 'Macroferia Regional de Talca'
$data[30,2] = 'Maule'
$data[30,3] = 44985
$data[30,4] = 7
$data[30,5] = 'Fruta'
$data[30,6] = 100107
$data[30,7] = 'Otros'
$data[30,8] = 100107011
$data[30,9] = 'Tuna'
$data[30,10] = 'Sin especificar'
$data[30,11] = 'Especial'
$data[30,12] = 300
$data[30,13] = 18000
$data[30,14] = 18000
$data[30,15] = 18000
$data[30,16] = '$/caja 18 kilos'
$data[30,17] = 'Provincia de Limarí'
$data[30,18] = 1000
$data[30,19] = 18

$data[31,0] = 5
$data[31,1] = 'Macroferia Regional de Talca'
$data[31,2] = 'Maule'
$data[31,3] = 44985
$data[31,4] = 7
$data[31,5] = 'Fruta'
$data[31,6] = 100107
$data[31,7] = 'Otros'
$data[31,8] = 100107011
$data[31,9] = 'Tuna'
$data[31,10] = 'Sin especificar'
$data[31,11] = 'Segunda'
$data[31,12] = 150
$data[31,13] = 12000
$data[31,14] = 12000
$data[31,15] = 12000
$data[31,16] = '$/caja 18 kilos'
$data[31,17] = 'Provincia de Limarí'
$data[31,18] = 667
$data[31,19] = 18

$data[32,0] = 5
$data[32,1] = 'Macroferia Regional de Talca'
$data[32,2] = 'Maule'
$data[32,3] = 44295
$data[32,4] = 7
$data[32,5] = 'Fruta'
$data[32,6] = 100107
$data[32,7] = 'Otros'
$data[32,8] = 100107011
$data[32,9] = 'Tuna'
$data[32,10] = 'Sin especificar'
$data[32,11] = 'Segunda'
$data[32,12] = 130
$data[32,13] = 10000
$data[32,14] = 10000
$data[32,15] = 10000
$data[32,16] = '$/caja 18 kilos'
$data[32,17] = 'Provincia de Melipilla'
$data[32,18] = 556
$data[32,19] = 18

$data[33,0] = 5
$data[33,1] = 'Macroferia Regional de Talca'
$data[33,2] = 'Maule'
$data[33,3] = 45083
$data[33,4] = 7
$data[33,5] = 'Fruta'
$data[33,6] = 100107
$data[33,7] = 'Otros'
$data[33,8] = 100107011
$data[33,9] = 'Tuna'
$data[33,10] = 'Sin especificar'
$data[33,11] = 'Primera'
$data[33,12] = 210
$data[33,13] = 20000
$data[33,14] = 20000
$data[33,15] = 20000
$data[33,16] = '$/caja 18 kilos'
$data[33,17] = 'Provincia de Melipilla'
$data[33,18] = 1111
$data[33,19] = 18

$data[34,0] = 5
$data[34,1] = 'Macroferia Regional de Talca'
$data[34,2] = 'Maule'
$data[34,3] = 45022
$data[34,4] = 7
$data[34,5] = 'Fruta'
$data[34,6] = 100107
$data[34,7] = 'Otros'
$data[34,8] = 100107011
$data[34,9] = 'Tuna'
$data[34,10] = 'Sin especificar'
$data[34,11] = 'Especial'
$data[34,12] = 200
$data[34,13] = 18000
$data[34,14] = 18000
$data[34,15] = 18000
$data[34,16] = '$/caja 18 kilos'
$data[34,17] = 'Provincia de Melipilla'
$data[34,18] = 1000
$data[34,19] = 18

$data[35,0] = 5
$data[35,1] = 'Macroferia Regional de Talca'
$data[35,2] = 'Maule'
$data[35,3] = 44271
$data[35,4] = 7
$data[35,5] = 'Fruta'
$data[35,6] = 100107
$data[35,7] = 'Otros'
$data[35,8] = 100107011
$data[35,9] = 'Tuna'
$data[35,10] = 'Sin especificar'
$data[35,11] = 'Primera'
$data[35,12] = 60
$data[35,13] = 15000
$data[35,14] = 15000
$data[35,15] = 15000
$data[35,16] = '$/caja 18 kilos'
$data[35,17] = 'Provincia de Melipilla'
$data[35,18] = 833
$data[35,19] = 18

$data[36,0] = 5
$data[36,1] = 'Macroferia Regional de Talca'
$data[36,2] = 'Maule'
$data[36,3] = 44350
$data[36,4] = 7
$data[36,5] = 'Fruta'
$data[36,6] = 100107
$data[36,7] = 'Otros'
$data[36,8] = 100107011
$data[36,9] = 'Tuna'
$data[36,10] = 'Sin especificar'
$data[36,11] = 'Especial'
$data[36,12] = 60
$data[36,13] = 24000
$data[36,14] = 24000
$data[36,15] = 24000
$data[36,16] = '$/caja 18 kilos'
$data[36,17] = 'Provincia de Limarí'
$data[36,18] = 1333
$data[36,19] = 18

$data[37,0] = 5
$data[37,1] = 'Macroferia Regional de Talca'
$data[37,2] = 'Maule'
$data[37,3] = 44224
$data[37,4] = 7
$data[37,5] = 'Fruta'
$data[37,6] = 100107
$data[37,7] = 'Otros'
$data[37,8] = 100107011
$data[37,9] = 'Tuna'
$data[37,10] = 'Sin especificar'
$data[37,11] = 'Primera'
$data[37,12] = 120
$data[37,13] = 18000
$data[37,14] = 18000
$data[37,15] = 18000
$data[37,16] = '$/caja 16 kilos'
$data[37,17] = 'Provincia de Limarí'
$data[37,18] = 1125
$data[37,19] = 16

$data[38,0] = 5
$data[38,1] = 'Macroferia Regional de Talca'
$data[38,2] = 'Maule'
$data[38,3] = 44298
$data[38,4] = 7
$data[38,5] = 'Fruta'
$data[38,6] = 100107
$data[38,7] = 'Otros'
$data[38,8] = 100107011
$data[38,9] = 'Tuna'
$data[38,10] = 'Sin especificar'
$data[38,11] = 'Extra (doble especial)'
$data[38,12] = 160
$data[38,13] = 20000
$data[38,14] = 20000
$data[38,15] = 20000
$data[38,16] = '$/caja 18 kilos'
$data[38,17] = 'Provincia de Melipilla'
$data[38,18] = 1111
$data[38,19] = 18

$data[39,0] = 5
$data[39,1] = 'Macroferia Regional de Talca'
$data[39,2] = 'Maule'
$data[39,3] = 44284
$data[39,4] = 7
$data[39,5] = 'Fruta'
$data[39,6] = 100107
$data[39,7] = 'Otros'
$data[39,8] = 100107011
$data[39,9] = 'Tuna'
$data[39,10] = 'Sin especificar'
$data[39,11] = 'Especial'
$data[39,12] = 120
$data[39,13] = 13000
$data[39,14] = 13000
$data[39,15] = 13000
$data[39,16] = '$/caja 18 kilos'
$data[39,17] = 'Provincia de Melipilla'
$data[39,18] = 722
$data[39,19] = 18

$data[40,0] = 5
$data[40,1] = 'Macroferia Regional de Talca'
$data[40,2] = 'Maule'
$data[40,3] = 44284
$data[40,4] = 7
$data[40,5] = 'Fruta'
$data[40,6] = 100107
$data[40,7] = 'Otros'
$data[40,8] = 100107011
$data[40,9] = 'Tuna'
$data[40,10] = 'Sin especificar'
$data[40,11] = 'Extra (doble especial)'
$data[40,12] = 100
$data[40,13] = 15000
$data[40,14] = 15000
$data[40,15] = 15000
$data[40,16] = '$/caja 18 kilos'
$data[40,17] = 'Provincia de Melipilla'
$data[40,18] = 833
$data[40,19] = 18

$data[41,0] = 5
$data[41,1] = 'Macroferia Regional de Talca'
$data[41,2] = 'Maule'
$data[41,3] = 44284
$data[41,4] = 7
$data[41,5] = 'Fruta'
$data[41,6] = 100107
$data[41,7] = 'Otros'
$data[41,8] = 100107011
$data[41,9] = 'Tuna'
$data[41,10] = 'Sin especificar'
$data[41,11] = 'Primera'
$data[41,12] = 50
$data[41,13] = 12000
$data[41,14] = 12000
$data[41,15] = 12000
$data[41,16] = '$/caja 18 kilos'
$data[41,17] = 'Provincia de Melipilla'
$data[41,18] = 667
$data[41,19] = 18

$data[42,0] = 5
$data[42,1] = 'Macroferia Regional de Talca'
$data[42,2] = 'Maule'
$data[42,3] = 44698
$data[42,4] = 7
$data[42,5] = 'Fruta'
$data[42,6] = 100107
$data[42,7] = 'Otros'
$data[42,8] = 100107011
$data[42,9] = 'Tuna'
$data[42,10] = 'Sin especificar'
$data[42,11] = 'Especial'
$data[42,12] = 150
$data[42,13] = 20000
$data[42,14] = 20000
$data[42,15] = 20000
$data[42,16] = '$/caja 18 kilos'
$data[42,17] = 'Provincia de Limarí'
$data[42,18] = 1111
$data[42,19] = 18

$data[43,0] = 5
$data[43,1] = 'Macroferia Regional de Talca'
$data[43,2] = 'Maule'
$data[43,3] = 44698
$data[43,4] = 7
$data[43,5] = 'Fruta'
$data[43,6] = 100107
$data[43,7] = 'Otros'
$data[43,8] = 100107011
$data[43,9] = 'Tuna'
$data[43,10] = 'Sin especificar'
$data[43,11] = 'Primera'
$data[43,12] = 180
$data[43,13] = 18000
$data[43,14] = 18000
$data[43,15] = 18000
$data[43,16] = '$/caja 18 kilos'
$data[43,17] = 'Provincia de Limarí'
$data[43,18] = 1000
$data[43,19] = 18

$data[44,0] = 5
$data[44,1] = 'Macroferia Regional de Talca'
$data[44,2] = 'Maule'
$data[44,3] = 45085
$data[44,4] = 7
$data[44,5] = 'Fruta'
$data[44,6] = 100107
$data[44,7] = 'Otros'
$data[44,8] = 100107011
$data[44,9] = 'Tuna'
$data[44,10] = 'Sin especificar'
$data[44,11] = 'Primera'
$data[44,12] = 150
$data[44,13] = 20000
$data[44,14] = 20000
$data[44,15] = 20000
$data[44,16] = '$/caja 18 kilos'
$data[44,17] = 'Provincia de Melipilla'
$data[44,18] = 1111
$data[44,19] = 18

$data[45,0] = 5
$data[45,1] = 'Macroferia Regional de Talca'
$data[45,2] = 'Maule'
$data[45,3] = 44301
$data[45,4] = 7
$data[45,5] = 'Fruta'
$data[45,6] = 100107
$data[45,7] = 'Otros'
$data[45,8] = 100107011
$data[45,9] = 'Tuna'
$data[45,10] = 'Sin especificar'
$data[45,11] = 'Primera'
$data[45,12] = 100
$data[45,13] = 16000
$data[45,14] = 16000
$data[45,15] = 16000
$data[45,16] = '$/caja 18 kilos'
$data[45,17] = 'Provincia de Melipilla'
$data[45,18] = 889
$data[45,19] = 18

$data[46,0] = 5
$data[46,1] = 'Macroferia Regional de Talca'
$data[46,2] = 'Maule'
$data[46,3] = 44267
$data[46,4] = 7
$data[46,5] = 'Fruta'
$data[46,6] = 100107
$data[46,7] = 'Otros'
$data[46,8] = 100107011
$data[46,9] = 'Tuna'
$data[46,10] = 'Sin especificar'
$data[46,11] = 'Primera'
$data[46,12] = 120
$data[46,13] = 13000
$data[46,14] = 13000
$data[46,15] = 13000
$data[46,16] = '$/caja 18 kilos'
$data[46,17] = 'Provincia de Melipilla'
$data[46,18] = 722
$data[46,19] = 18

$ws.Range("A2:T48").Value = $data
